$d = $word.ActiveDocument

# Locate the "Repository" heading paragraph (exact, whole-word match).
$range = $d.Content
$found = $range.Find.Execute("Repository", $true, $true, $false, $false, `
                              $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to a zero-length range right after the word "Repository"
    # (i.e. right before its paragraph mark) and insert a brand-new
    # paragraph there via raw OOXML so it does not inherit the Heading2
    # style / numbering that the "Repository" paragraph carries.
    $insertionPoint = $d.Range($range.End, $range.End)

    $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:ind w:left="708"/>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
</w:rPr>
<w:t xml:space="preserve">Project added to </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
</w:rPr>
<w:t>https://github.com/mustafakoroglu/messageQueue/</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

    $insertionPoint.InsertXML($xml) | Out-Null
}
